$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "Num"
$ws.Range("B23").Value = "Top"
$ws.Range("C23").Value = "Bottom"
$ws.Range("D23").Value = "Difference"

$ws.Range("A24").Value = 1
$ws.Range("B24").Value = 8.402608871459961
$ws.Range("C24").Value = 0.6628193855285645
$ws.Range("D24").Value = 7.739789485931396

$ws.Range("A25").Value = 2
$ws.Range("B25").Value = 8.299276351928711
$ws.Range("C25").Value = 0.7217079997062683
$ws.Range("D25").Value = 7.577568531036377

$ws.Range("A26").Value = 3
$ws.Range("B26").Value = 8.402608871459961
$ws.Range("C26").Value = 0.6619071364402771
$ws.Range("D26").Value = 7.740701675415039

$ws.Range("A28").Value = "Num"
$ws.Range("B28").Value = "Top"
$ws.Range("C28").Value = "Bottom"
$ws.Range("D28").Value = "Difference"

$ws.Range("A29").Value = 1
$ws.Range("B29").Value = 8.402608871459961
$ws.Range("C29").Value = 0.6563481092453003
$ws.Range("D29").Value = 7.746260643005371

$ws.Range("A30").Value = 2
$ws.Range("B30").Value = 8.299276351928711
$ws.Range("C30").Value = 0.7046168446540833
$ws.Range("D30").Value = 7.594659328460693

$ws.Range("A31").Value = 3
$ws.Range("B31").Value = 8.402608871459961
$ws.Range("C31").Value = 0.6621693968772888
$ws.Range("D31").Value = 7.740439414978027

$ws.Range("A33").Value = "Num"
$ws.Range("B33").Value = "Top"
$ws.Range("C33").Value = "Bottom"
$ws.Range("D33").Value = "Difference"

$ws.Range("A34").Value = 1
$ws.Range("B34").Value = 8.402608871459961
$ws.Range("C34").Value = 0.6620481610298157
$ws.Range("D34").Value = 7.740560531616211

$ws.Range("A35").Value = 2
$ws.Range("B35").Value = 8.299276351928711
$ws.Range("C35").Value = 0.7242947220802307
$ws.Range("D35").Value = 7.574981689453125

$ws.Range("A37").Value = "Num"
$ws.Range("B37").Value = "Top"
$ws.Range("C37").Value = "Bottom"
$ws.Range("D37").Value = "Difference"

$ws.Range("A38").Value = 1
$ws.Range("B38").Value = 8.402608871459961
$ws.Range("C38").Value = 0.6635361313819885
$ws.Range("D38").Value = 7.739072799682617

$ws.Range("A39").Value = 6
$ws.Range("B39").Value = 3.814814567565918
$ws.Range("C39").Value = 0.9300078749656677
$ws.Range("D39").Value = 2.884806632995605

$ws.Range("A40").Value = 13
$ws.Range("B40").Value = 8.402608871459961
$ws.Range("C40").Value = 0.6647151708602905
$ws.Range("D40").Value = 7.737893581390381

$ws.Range("A42").Value = "Num"
$ws.Range("B42").Value = "Top"
$ws.Range("C42").Value = "Bottom"
$ws.Range("D42").Value = "Difference"

$ws.Range("A43").Value = 1
$ws.Range("B43").Value = 8.402608871459961
$ws.Range("C43").Value = 0.6621693968772888
$ws.Range("D43").Value = 7.740439414978027

$ws.Range("A45").Value = "Num"
$ws.Range("B45").Value = "Top"
$ws.Range("C45").Value = "Bottom"
$ws.Range("D45").Value = "Difference"

$ws.Range("A46").Value = 1
$ws.Range("B46").Value = 8.402608871459961
$ws.Range("C46").Value = 0.6811108589172363
$ws.Range("D46").Value = 7.721498012542725

$ws.Range("A48").Value = "Num"
$ws.Range("B48").Value = "Top"
$ws.Range("C48").Value = "Bottom"
$ws.Range("D48").Value = "Difference"

$ws.Range("A49").Value = 1
$ws.Range("B49").Value = 8.402608871459961
$ws.Range("C49").Value = 0.6569118499755859
$ws.Range("D49").Value = 7.745697021484375

$ws.Range("A51").Value = 17
$ws.Range("B51").Value = 8.402608871459961
$ws.Range("C51").Value = 0.6690974235534668
$ws.Range("D51").Value = 7.733511447906494

$ws.Range("A53").Value = "Num"
$ws.Range("B53").Value = "Top"
$ws.Range("C53").Value = "Bottom"
$ws.Range("D53").Value = "Difference"

$ws.Range("A54").Value = 1
$ws.Range("B54").Value = 8.402608871459961
$ws.Range("C54").Value = 0.6624118089675903
$ws.Range("D54").Value = 7.74019718170166

$ws.Range("A55").Value = 2
$ws.Range("B55").Value = 8.299276351928711
$ws.Range("C55").Value = 0.7033462524414062
$ws.Range("D55").Value = 7.595930099487305

$ws.Range("A56").Value = 3
$ws.Range("B56").Value = 7.370575428009033
$ws.Range("C56").Value = 0.7949501276016235
$ws.Range("D56").Value = 6.575625419616699

$ws.Range("A57").Value = 4
$ws.Range("B57").Value = 7.817458152770996
$ws.Range("C57").Value = 0.6562545299530029
$ws.Range("D57").Value = 7.161203384399414

$ws.Range("A58").Value = 5
$ws.Range("B58").Value = 7.809852600097656
$ws.Range("C58").Value = 0.9485856294631958
$ws.Range("D58").Value = 6.86126708984375

$ws.Range("A59").Value = 4
$ws.Range("B59").Value = 8.402608871459961
$ws.Range("C59").Value = 0.6605494022369385
$ws.Range("D59").Value = 7.742059707641602

$ws.Range("A61").Value = 5
$ws.Range("B61").Value = 8.299276351928711
$ws.Range("C61").Value = 0.6904088854789734
$ws.Range("D61").Value = 7.608867645263672

# Copy formatting (bold, border, centered/top alignment) from the existing header row (row 3) onto the new header rows
$ws.Range("A3:D3").Copy()
$ws.Range("A23:D23").PasteSpecial(-4122)
$ws.Range("A28:D28").PasteSpecial(-4122)
$ws.Range("A33:D33").PasteSpecial(-4122)
$ws.Range("A37:D37").PasteSpecial(-4122)
$ws.Range("A42:D42").PasteSpecial(-4122)
$ws.Range("A45:D45").PasteSpecial(-4122)
$ws.Range("A48:D48").PasteSpecial(-4122)
$ws.Range("A53:D53").PasteSpecial(-4122)
$excel.CutCopyMode = 0
